$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.293.85"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "2.958.66"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.90"
$ws.Range("E5").Value = "  +2.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.83"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.541"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.593"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.04"
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0841"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "3.440.13"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.41"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.47"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "2.977.58"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.963"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "51.357.52"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.34"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.90"
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.93"
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.96"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  +4.26%  "
$ws.Range("B26").Value = "Filecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.50"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.169"
$ws.Range("E27").Value = "  -3.54%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.16"
$ws.Range("E28").Value = "  +16.61%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.89"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.112"
$ws.Range("E31").Value = "  +7.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.82"
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.67"
$ws.Range("E33").Value = "  -1.99%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.08"
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.02"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0447"
$ws.Range("E36").Value = "  +5.04%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.07"
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.29"
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.58"
$ws.Range("E40").Value = "  -4.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.84"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "123.43"
$ws.Range("E43").Value = "  +3.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.12"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("E45").Value = "  +18.96%  "
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.35"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("D48").Value = "2.034.49"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.22"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0352"
$ws.Range("E50").Value = "  +10.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.15"
$ws.Range("E51").Value = "  +2.07%  "
